$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$defaultStyle = $ws.Range("D4").Style

$ws.Range("D2").Value = '58.944.53'
$ws.Range("E2").Value = '  +4.41%  '
$ws.Range("D3").Value = '3.304.70'
$ws.Range("E3").Value = '  +1.57%  '
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = "'408.18"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = '  +2.32%  '
$ws.Range("D6").Value = "'111.91"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = '  +0.80%  '
$ws.Range("D7").Value = "'0.583"
$ws.Range("D7").Style = $defaultStyle
$ws.Range("E7").Value = '  +4.32%  '
$ws.Range("E8").Value = '  +0.23%  '
$ws.Range("D9").Value = "'0.629"
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Value = '  +1.57%  '
$ws.Range("D10").Value = "'39.78"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = '  +0.78%  '
$ws.Range("D11").Value = "'0.0977"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = '  +3.48%  '
$ws.Range("E12").Value = '  +1.26%  '
$ws.Range("D13").Value = '3.850.25'
$ws.Range("E13").Value = '  +2.24%  '
$ws.Range("D14").Value = "'8.42"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = '  +3.93%  '
$ws.Range("D15").Value = "'19.31"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = '  +0.43%  '
$ws.Range("D16").Value = '3.288.17'
$ws.Range("E16").Value = '  +1.19%  '
$ws.Range("D17").Value = "'1.03"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = '  -1.13%  '
$ws.Range("D18").Value = '58.941.32'
$ws.Range("E18").Value = '  +4.55%  '
$ws.Range("D19").Value = "'10.65"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = '  -3.27%  '
$ws.Range("D20").Value = "'3.32"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("D21").Value = "'0.0000109"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = '  +5.80%  '
$ws.Range("D22").Value = "'13.03"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = '  -0.24%  '
$ws.Range("D23").Value = "'301.93"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = '  +0.87%  '
$ws.Range("D24").Value = "'74.98"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("E25").Value = '  -0.35%  '
$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").Value = "'4.47"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = '  +2.49%  '
$ws.Range("D27").Value = "'28.34"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = '  +0.64%  '
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").Value = "'0.179"
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").Value = '  +5.74%  '
$ws.Range("D29").Value = "'7.81"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = '  -3.48%  '
$ws.Range("D30").Value = "'7.52"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = '  +2.58%  '
$ws.Range("D31").Value = "'0.115"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = '  +3.14%  '
$ws.Range("D32").Value = "'1.00"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("D33").Value = "'11.44"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = '  +3.34%  '
$ws.Range("D34").Value = "'39.32"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = '  +2.00%  '
$ws.Range("D35").Value = "'0.0513"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = '  +5.25%  '
$ws.Range("D36").Value = "'51.98"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = '  +0.84%  '
$ws.Range("E37").Value = '  -1.86%  '
$ws.Range("D38").Value = "'3.10"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = '  -2.32%  '
$ws.Range("E39").Value = '  +0.15%  '
$ws.Range("D40").Value = "'3.37"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = '  -4.21%  '
$ws.Range("D41").Value = "'137.83"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = '  +2.95%  '
$ws.Range("D42").Value = "'0.122"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = '  +1.59%  '
$ws.Range("D43").Value = "'1.90"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = '  -2.14%  '
$ws.Range("D44").Value = "'16.82"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = '  -4.28%  '
$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").Value = "'0.281"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = '  -1.27%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").Value = "'3.90"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = '  -1.96%  '
$ws.Range("D47").Value = "'2.28"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = '  +8.99%  '
$ws.Range("D48").Value = "'22.44"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = '  +1.34%  '
$ws.Range("D49").Value = '2.203.65'
$ws.Range("E49").Value = '  +2.54%  '
$ws.Range("E50").Value = '  -0.07%  '
$ws.Range("E51").Value = '  -5.05%  '
